$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Sheet6 -> Sheet3)
$ws.Name = "Sheet3"

# Update the report generation date/time (row 1)
$ws.Range("D1").Value = 45572
$ws.Range("F1").Value = 0.808196076388889

# Rename the "Induction Hardening Bearing Surface 1" process to
# "Induction Hardening Bearing Surfaces 1, 2" (less impactful hardening process)
$ws.Range("W16").Value = "Induction Hardening Bearing Surfaces 1, 2"
$ws.Range("B19").Value = "Induction Hardening Bearing Surfaces 1, 2"

# Update the impact values for the (renamed) hardening process row
$ws.Range("E19").Value = 252
$ws.Range("W19").Value = 252
